$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells J1:N1 (plain numbers 1-5) and O1 ("Gemiddelde", bold) ---
$ws.Range("J1").Value = 1
$ws.Range("K1").Value = 2
$ws.Range("L1").Value = 3
$ws.Range("M1").Value = 4
$ws.Range("N1").Value = 5
$ws.Range("O1").Value = "Gemiddelde"
$ws.Range("O1").Font.Bold = $true

# --- Column O width ---
$ws.Range("O1").ColumnWidth = 14.11

# --- H column: convert displayed percentages to fractions (existing values /10) ---
$ws.Range("H2").Value = 0.028
$ws.Range("H3").Value = 0.067
$ws.Range("H11").Value = 0.016
$ws.Range("H12").Value = 0.022
$ws.Range("H13").Value = 0.022
$ws.Range("H14").Value = 0.018
$ws.Range("H15").Value = 0.017

# --- Row 12 new measurement data + average ---
$ws.Range("J12").Value = 11897.3
$ws.Range("K12").Value = 11522.7
$ws.Range("L12").Value = 13794
$ws.Range("M12").Value = 13949
$ws.Range("N12").Value = 14138.6
$ws.Range("O12").Formula = "=AVERAGE(J12:N12)"

# --- Row 13 new measurement data + average ---
$ws.Range("J13").Value = 13026.8
$ws.Range("K13").Value = 13889.3
$ws.Range("L13").Value = 13740.1
$ws.Range("M13").Value = 14030.9
$ws.Range("N13").Value = 14141
$ws.Range("O13").Formula = "=AVERAGE(J13:N13)"

# --- Row 14 new measurement data + average ---
$ws.Range("J14").Value = 10965.8
$ws.Range("K14").Value = 13431.8
$ws.Range("L14").Value = 14317
$ws.Range("M14").Value = 13389.4
$ws.Range("N14").Value = 13262.3
$ws.Range("O14").Formula = "=AVERAGE(J14:N14)"

# --- Row 15 new measurement data + average ---
$ws.Range("J15").Value = 11069.8
$ws.Range("K15").Value = 11589.3
$ws.Range("L15").Value = 10944
$ws.Range("M15").Value = 11208
$ws.Range("N15").Value = 11019.8
$ws.Range("O15").Formula = "=AVERAGE(J15:N15)"

# --- View: select the freshly added block, active cell at the bottom-right (O15) ---
$ws.Range("J12:O15").Select()
